$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A47").Value = "SP25092022112318"
$ws.Range("B47").Value = 100
$ws.Range("C47").Value = "PANADOL STRIP 10"
$ws.Range("D47").Value = 1
$ws.Range("E47").Value = 15
